# "Inscription et nouveau logo"
# Adds new "revenue" rows (E3:F6) for rallye/CIME registration income on
# Sheet1, renumbering the shared-string table as a side effect, and
# switches the active sheet/selection from Sheet1 to JeunesCotisation,
# along with a handful of column-width tweaks on both sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("JeunesCotisation")

# --- New income rows on Sheet1 (E3:F6) -------------------------------
$ws1.Range("E3").Value = "Rallye Castor 2016"
$ws1.Range("F3").Value = 350

$ws1.Range("E4").Value = "Rallye Branche Jaune"
$ws1.Range("F4").Value = 200

$ws1.Range("E5").Value = "CIME"
$ws1.Range("F5").Value = 250

$ws1.Range("E6").Value = "Budget 10 ans"

# --- Column width tweaks ----------------------------------------------
# Sheet1: columns A, E, L get new widths
$ws1.Columns.Item(1).ColumnWidth = 13.1666666666667
$ws1.Columns.Item(5).ColumnWidth = 19.5
$ws1.Columns.Item(12).ColumnWidth = 24.5

# JeunesCotisation: columns A, B, C, D get new widths
$ws2.Columns.Item(1).ColumnWidth = 18.8333333333333
$ws2.Columns.Item(2).ColumnWidth = 10.5
$ws2.Columns.Item(3).ColumnWidth = 15.0
$ws2.Columns.Item(4).ColumnWidth = 10.5

# --- Switch the active sheet / selection -------------------------------
# Sheet1 ends up unselected, with the cursor left on E7
[void]$ws1.Range("E7").Select()

# JeunesCotisation becomes the active tab, with the cursor on A3
[void]$ws2.Activate()
[void]$ws2.Range("A3").Select()
